$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert 8 blank rows at position 10 (below the existing block start),
# which pushes the current rows 10-21 down to 18-29 without touching the
# header's formatting boundary (avoids style bleed at row 2).
$ws.Rows("10:17").Insert()

# Step 2: move the data currently sitting in rows 2:9 down into the newly
# freed rows 10:17 (this is the original first 8 data rows).
$ws.Range("A2:C9").Cut($ws.Range("A10:C17"))

# Step 3: populate rows 2:9 with the new data that belongs at the top.
$newTop = @(
    @(0.00580321977447178, 0.07317293116024544, 0.01743147575429496),
    @(-0.01212380492907687, 0.1479228914392232, 0.03119464347861243),
    @(-0.05563860289676467, 0.2476778115545005, 0.03135671048444155),
    @(0.000710598117082524, -0.0459115079471038, -0.03720979673825954),
    @(0.03395912094915998, -0.002209710983597735, -0.02471199281969851),
    @(0.0100699262883591, -0.003861541194575109, -0.02743906991518268),
    @(-0.04216528505238931, -0.05587235412427344, -0.005946585338334982),
    @(-0.2539235970803638, -0.6386596262454889, 0.1363320350646959)
)

$r = 2
foreach ($row in $newTop) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Step 4: append 2 brand-new rows at the bottom (rows 30 and 31).
$newBottom = @(
    @(0.5868015289306701, 4.023616756711703, 0.7568838426044971),
    @(-0.3439888250538894, 1.417871174155451, 1.152574896812441)
)

$r = 30
foreach ($row in $newBottom) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
